$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("commondata")

# Update the URL value held in B1 (keeps its existing hyperlink + style)
$ws.Range("B1").Value = "http://192.168.235.128:8080/crm/ShowHomePage.do"

# Move the active selection to B1 (matches the saved selection in the edit)
$null = $ws.Range("B1").Select()
